# Apply updated crypto price / volume data (Wed Jun  5 03:45:53 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.812.41'
$ws.Range("E2").Value = '  +2.32%  '
$ws.Range("D3").Value = '3.809.06'
$ws.Range("E3").Value = '  +0.81%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '700.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +11.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.20%  '
$ws.Range("D7").Value = '3.807.81'
$ws.Range("E7").Value = '  +0.79%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +0.92%  '
$ws.Range("E10").Value = '  +2.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.66'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +12.86%  '
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("E13").Value = '  +4.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.15'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.68%  '
$ws.Range("D15").Value = '4.450.27'
$ws.Range("D16").Value = '3.808.48'
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("D17").Value = '70.874.33'
$ws.Range("E17").Value = '  +2.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.71'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.21'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.32'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +18.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '479.97'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.45%  '
$ws.Range("E23").Value = '  +1.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.81'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.09%  '
$ws.Range("E25").Value = '  +0.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.28%  '
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("E28").Value = '  +3.03%  '
$ws.Range("D29").Value = '3.960.20'
$ws.Range("E29").Value = '  +0.87%  '
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.10'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +15.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.31'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.51'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.55'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.178'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.20'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.35%  '
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("D38").Value = '3.759.00'
$ws.Range("E38").Value = '  +0.80%  '
$ws.Range("E39").Value = '  +1.33%  '
$ws.Range("E40").Value = '  +6.94%  '
$ws.Range("E41").Value = '  +3.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.000334'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +24.07%  '
$ws.Range("E43").Value = '  +12.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.968'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.42%  '
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '45.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '159.93'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '48.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.69%  '
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.299'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.37%  '
